$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A9").Value = "Spring"
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = "http://www.mcmaster.com/#9657k61/=x839tg"

$ws.Range("A9").Select()
